$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B
$ws.Columns("B:B").Insert()

# Set the header text for the new column B (row 1)
$ws.Range("B1").Value = "比賽年份 Year of Competition"

# Update selection to C3 as seen in final file
$ws.Range("C3").Select()
